$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.036.65'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').Value = '2.348.97'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('E4').Value = '  +0.05%  '
$c = $ws.Range('D5')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.677'
$c.Style = $s

$ws.Range('E5').Value = '  +0.02%  '
$c = $ws.Range('D6')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '239.03'
$c.Style = $s

$ws.Range('E6').Value = '  +1.63%  '
$c = $ws.Range('D7')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '73.96'
$c.Style = $s

$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('E8').Value = '  +0.00%  '
$c = $ws.Range('D9')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.594'
$c.Style = $s

$ws.Range('E9').Value = '  +8.95%  '
$c = $ws.Range('D10')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.100'
$c.Style = $s

$ws.Range('E10').Value = '  +1.95%  '
$ws.Range('E11').Value = '  +0.18%  '
$c = $ws.Range('D12')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '31.99'
$c.Style = $s

$ws.Range('E12').Value = '  +13.29%  '
$c = $ws.Range('D13')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.107'
$c.Style = $s

$ws.Range('E13').Value = '  +0.73%  '
$c = $ws.Range('D14')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '7.17'
$c.Style = $s

$ws.Range('E14').Value = '  +7.01%  '
$ws.Range('D15').Value = '2.698.18'
$ws.Range('E15').Value = '  -0.12%  '
$c = $ws.Range('D16')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '16.49'
$c.Style = $s

$ws.Range('E16').Value = '  -1.24%  '
$c = $ws.Range('D17')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.897'
$c.Style = $s

$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').Value = '2.328.87'
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('D19').Value = '43.925.18'
$ws.Range('E19').Value = '  +0.51%  '
$c = $ws.Range('D20')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0000102'
$c.Style = $s

$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('E21').Value = '  +3.49%  '
$c = $ws.Range('D22')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '76.77'
$c.Style = $s

$ws.Range('E22').Value = '  -0.56%  '
$c = $ws.Range('D23')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '256.00'
$c.Style = $s

$ws.Range('E23').Value = '  +1.18%  '
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D24')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.91'
$c.Style = $s

$ws.Range('E24').Value = '  +19.75%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D25')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = $s

$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('E26').Value = '  -1.30%  '
$c = $ws.Range('D27')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.49'
$c.Style = $s

$ws.Range('E27').Value = '  +0.09%  '
$c = $ws.Range('D28')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '10.65'
$c.Style = $s

$ws.Range('E28').Value = '  +0.62%  '
$c = $ws.Range('D29')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.25'
$c.Style = $s

$ws.Range('E29').Value = '  -0.96%  '
$c = $ws.Range('D30')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '22.70'
$c.Style = $s

$ws.Range('E30').Value = '  +1.63%  '
$c = $ws.Range('D31')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '174.61'
$c.Style = $s

$ws.Range('E31').Value = '  +1.42%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D32')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.137'
$c.Style = $s

$ws.Range('E32').Value = '  +3.61%  '
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D33')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.126'
$c.Style = $s

$ws.Range('E33').Value = '  -2.71%  '
$ws.Range('E34').Value = '  +5.65%  '
$c = $ws.Range('D35')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '5.27'
$c.Style = $s

$ws.Range('E35').Value = '  +1.88%  '
$ws.Range('E36').Value = '  +3.77%  '
$c = $ws.Range('D37')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.71'
$c.Style = $s

$ws.Range('E37').Value = '  -3.83%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D38')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.35'
$c.Style = $s

$ws.Range('E38').Value = '  -2.47%  '
$ws.Range('B39').Value = 'THORChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$c = $ws.Range('D39')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '6.32'
$c.Style = $s

$ws.Range('E39').Value = '  -1.06%  '
$c = $ws.Range('D40')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0281'
$c.Style = $s

$ws.Range('E40').Value = '  +4.55%  '
$c = $ws.Range('D41')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.108'
$c.Style = $s

$ws.Range('E41').Value = '  +11.33%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D42')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '19.03'
$c.Style = $s

$ws.Range('E42').Value = '  -2.15%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D43')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '9.04'
$c.Style = $s

$ws.Range('E43').Value = '  +1.99%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D44')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.201'
$c.Style = $s

$ws.Range('E44').Value = '  +8.64%  '
$ws.Range('E45').Value = '  +0.00%  '
$c = $ws.Range('D46')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '4.66'
$c.Style = $s

$ws.Range('E46').Value = '  +4.76%  '
$ws.Range('B47').Value = 'MultiversX'
$ws.Range('C47').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$c = $ws.Range('D47')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '57.64'
$c.Style = $s

$ws.Range('E47').Value = '  +10.42%  '
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D48')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.25'
$c.Style = $s

$ws.Range('E48').Value = '  +1.67%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D49')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.47'
$c.Style = $s

$ws.Range('E49').Value = '  +8.19%  '
$ws.Range('E50').Value = '  +0.85%  '
$c = $ws.Range('D51')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '99.80'
$c.Style = $s

$ws.Range('E51').Value = '  +2.55%  '